# React_AW1.pptx - slide 3 ("SEO" panel): fix the subtitle text
# ("Search Engine Optimizaon" -> "Search Engine Optimization") and widen
# the placeholder (flipH="1") so the longer text still fits, shifting its
# left edge so the box keeps the same right edge.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes("Google Shape;400;p34")

# New off/ext (EMU, from the target OOXML):
#   off  x="4654018" y="3585925"  (y unchanged)
#   ext cx="2224057" cy="413100"  (cy unchanged)
# PowerPoint's COM surface works in points (1 pt = 12700 EMU); the literals
# below are chosen so the round-tripped EMU values land exactly on target.
$shp.Left = 366.4581298828125
$shp.Width = 175.1226043701172

# Replace the two runs ("Search Engine " + "Optimizaon") with a single,
# correctly spelled run ("Search Engine Optimization"). Re-selecting the
# whole range via Characters() before assigning Text keeps the result as
# one clean run instead of leaving behind stray spell-check run splits.
$tr = $shp.TextFrame.TextRange
$full = $tr.Characters(1, $tr.Length)
$full.Text = "Search Engine Optimization"
